$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C2").Value = "[CodePen](https://codepen.io/maiCoding/pen/XoaQKG)"
$ws.Range("C2").Select()
